# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 207 (pushing the existing
# rows 207-218 down to 208-219) in the Mango price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 207 - this shifts the old
# rows 207..218 down to 208..219 and copies formatting from the row above.
$ws.Rows.Item(207).Insert()

$row = 207

$ws.Cells.Item($row, 1).Value  = 7
$ws.Cells.Item($row, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value  = "Ñuble"
$ws.Cells.Item($row, 4).Value  = 45267
$ws.Cells.Item($row, 5).Value  = 16
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100108
$ws.Cells.Item($row, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value  = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 140
$ws.Cells.Item($row, 14).Value = 12000
$ws.Cells.Item($row, 15).Value = 13000
$ws.Cells.Item($row, 16).Value = 12571
$ws.Cells.Item($row, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item($row, 18).Value = "Perú"
$ws.Cells.Item($row, 19).Value = 3143
$ws.Cells.Item($row, 20).Value = 4
